$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (serial 45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = (Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0)

# Update prices in column D for rows 34-39
$ws.Range("D34").Value = 368
$ws.Range("D35").Value = 411
$ws.Range("D36").Value = 457
$ws.Range("D37").Value = 411
$ws.Range("D38").Value = 445
$ws.Range("D39").Value = 505
